$d = $word.ActiveDocument

# 1. Add a first-line indent (0.5in / 720 twips / 36pt) to the first paragraph
#    ("CIS 21JA Assignment 2 ... Name: Tom Ekshtein").
$d.Paragraphs(1).Range.ParagraphFormat.FirstLineIndent = 36

# 2. Question 2 (pipelined) answer: k + (2n-1)/24 cycles -> k + (n-1)/14 cycles
$d.Content.Find.Execute(
    "For k stages and n instructions, the number of required cycles is: k + (2n-1). Therefore, 5 + (2*10-1) = 24 cycles.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "For k stages and n instructions, the number of required cycles is: k + (n-1). Therefore, 5 + (10-1) = 14 cycles.",
    2)

# 3. Question 2 (non-pipelined) answer: rewritten to mention "2 cycles per stage"
#    and updated math (k*2n, 5*10(2)=100 cycles).
$d.Content.Find.Execute(
    "A non-pipelined processor, the number of required cycles with n instructions and k stages is: k*n. Therefore, 5*10=50 cycles ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "A non-pipelined processor, the number of required cycles with n instructions and k stages is where each stage takes 2 cycles is: k*2n. Therefore, 5*10(2)=100 cycles ",
    2)

# 4. Question 4 answer: "2 GB" -> full explanation of 2^32 / 2^16 address space sizes
$d.Content.Find.Execute(
    "2 GB",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Since 2^32bits = 4,294,967,296 unique addresses, or 4GB of memory. Therefore, 2^16bits = 65536 unique addresses or 64KB of memory.",
    2)
